$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Replace the long "Objetivos" essay text (B10/C10) with the professor id/name string.
$ws.Range("B10").Value2 = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C10").Value2 = "8554681 - Pedro Felipe Arce Castillo"

# 2) Delete row 13 entirely (it held only the professor name under "Docentes responsaveis:",
#    with no label in column A). This shifts rows 14-24 up to become rows 13-23, carrying
#    their row heights and styles with them.
$ws.Rows.Item(13).Delete()

# 3) Fix up the cells whose content needs to differ from what a pure shift would produce.

# New row 13 (was old row 14) - "Programa resumido:" summary becomes "Semestral"
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"

# New row 15 (was old row 16) - "Programa:" full syllabus becomes the activation date string
$ws.Range("B15").Value2 = "01/01/2012"
$ws.Range("C15").Value2 = "01/01/2012"

# New row 18 (was old row 19, "Metodo:") now needs the professor id/name string
$ws.Range("B18").Value2 = "8554681 - Pedro Felipe Arce Castillo"
$ws.Range("C18").Value2 = "8554681 - Pedro Felipe Arce Castillo"

# New row 19 (was old row 20, "Criterio:") now needs the evaluation-method text
$ws.Range("B19").Value2 = "A avaliação será feita por meio de duas provas escritas (P1 e P2)."
$ws.Range("C19").Value2 = "A avaliação será feita por meio de duas provas escritas (P1 e P2)."

# New row 20 (was old row 21, "Norma de recuperação:") now needs the final-grade formula text
$ws.Range("B20").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2"
$ws.Range("C20").Value2 = "A Nota final (NF) será calculada da seguinte maneira: NF = (P1 + P2)/2"

# New row 21 (was old row 22, "Bibliografia:") now needs the recovery-exam formula text
$ws.Range("B21").Value2 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value2 = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

# New row 23 (was old row 24) - prerequisite text, trailing newline trimmed
$ws.Range("B23").Value2 = "LOQ4087 -  Termodinâmica Química Aplicada I  (Requisito fraco)"
$ws.Range("C23").Value2 = "LOQ4087 -  Termodinâmica Química Aplicada I  (Requisito fraco)"
